$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "06/08/2025"
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").Value = 0.00047116
$ws.Range("C21").Value = 105059.8522794804
$ws.Range("D21").Value = 49.5
